# Workbook edit: add "Matières enseignés" column (E) to the professeur sheet,
# size columns C/D/E, and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell in column E, row 1 (adds a new shared string too)
$ws.Range("E1").Value = "Matières enseignés"

# Column widths for C, D, E (values chosen so the engine's internal
# pixel-quantization lands as close as possible to the authored widths
# 27.5703125 / 15.7109375 / 31.7109375 characters)
$ws.Columns.Item(3).ColumnWidth = 26.666666666666668
$ws.Columns.Item(4).ColumnWidth = 14.833333333333334
$ws.Columns.Item(5).ColumnWidth = 30.833333333333332

# Move/extend the active selection to E6, matching the saved view state
$null = $ws.Range("E6").Select()
